# Generate Report for Handback
# Applies the "handback" report-generation edits to the localization-status
# workbook: updates the status text, stamps the handback datetime, fills in
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for zh-cn and de-de, widens a couple of columns that now
# hold longer text, and relinks the per-row hyperlinks.

$wb = $excel.ActiveWorkbook

$srcCommit = "b8201a664a09ba5e0fda584516d18e680f0ccbd8"
$repoBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e"

# ---------------------------------------------------------------------
# 1) Overview sheet: widen the zh-cn / de-de status columns (E, F) now
#    that the status text is longer ("Handed back: in sync with en-US").
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# 2) Per-language sheets (zh-cn, de-de): update Status, Latest Target
#    File, Latest Handback File and Latest Handback DateTime, and widen
#    the Status (C) / Latest Handback File (J) columns to fit.
# ---------------------------------------------------------------------
$langs = @(
    @{ Sheet = "zh-cn"; XlfTail = "zh-cn.xlf"; HandbackTime = "2016-09-05 08:50:33" },
    @{ Sheet = "de-de"; XlfTail = "de-de.xlf"; HandbackTime = "2016-09-05 08:50:41" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(10).ColumnWidth = 39.17

    $handbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2." + $lang.XlfTail

    for ($row = 2; $row -le 3; $row++) {
        $ws.Cells.Item($row, 3).Value = "Handed back: in sync with en-US"   # C: Status
        $ws.Cells.Item($row, 9).Value = "a.md"                               # I: Latest Target File
        $ws.Cells.Item($row, 10).Value = $handbackFile                       # J: Latest Handback File
        $ws.Cells.Item($row, 11).Value = $lang.HandbackTime                  # K: Latest Handback DateTime
    }

    # Rebuild the hyperlinks for this sheet so A2/A3 (source file) and
    # I2/I3 (target file) all point at the right place, in row order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 1), "$repoBase/a.md", [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 9), "$repoBase/a.md", [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "$repoBase/b.md", [Type]::Missing, [Type]::Missing, "b.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 9), "$repoBase/a.md", [Type]::Missing, [Type]::Missing, "a.md")
}
